{"js": "// Update the \"contact e-mail\" paragraph in the poster template:\n//  - pluralise \"e-mail address\" -> \"e-mail addresses\"\n//  - add the new \"m.xochicale@ucl.ac.uk,\" address before the existing\n//    \"s.chopra@ucl.ac.uk\" one\n//  - add a trailing space after the last address (as in the target diff)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that still has the single-recipient sentence.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"Please contact the following e-mail address\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the 'Please contact the following e-mail address\u2026' paragraph.\");\n}\n\n// 1) \"address\" -> \"addresses\" (insert \"es\" right after the existing word).\nconst addressHits = target.search(\"address\", { matchCase: true, matchWholeWord: false });\naddressHits.load(\"items\");\nawait context.sync();\nif (addressHits.items.length === 0) {\n  throw new Error(\"Could not find 'address' to pluralise.\");\n}\naddressHits.items[0].insertText(\"es\", \"After\");\nawait context.sync();\n\n// 2) Insert the new e-mail address immediately before the existing one.\nconst emailHits = target.search(\"s.chopra@ucl.ac.uk\", { matchCase: true });\nemailHits.load(\"items\");\nawait context.sync();\nif (emailHits.items.length === 0) {\n  throw new Error(\"Could not find 's.chopra@ucl.ac.uk' to anchor the new address on.\");\n}\nemailHits.items[0].insertText(\"m.xochicale@ucl.ac.uk, \", \"Before\");\nawait context.sync();\n\n// 3) Add the trailing space that now follows the last e-mail address.\ntarget.getRange(\"End\").insertText(\" \", \"End\");\nawait context.sync();\n", "ps1": "# Update the \"contact e-mail\" paragraph in the poster template:\n#  - pluralise \"e-mail address\" -> \"e-mail addresses\"\n#  - add the new \"m.xochicale@ucl.ac.uk,\" address before the existing\n#    \"s.chopra@ucl.ac.uk\" one\n#  - add a trailing space after the last address (as in the target diff)\n\n$d = $word.ActiveDocument\n\n# 1) \"address\" -> \"addresses\": find the lead-in sentence and insert \"es\"\n#    right after the word \"address\".\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found1 = $rng.Find.Execute(\"Please contact the following e-mail address\")\nif (-not $found1) {\n    throw \"Could not find the 'Please contact the following e-mail address\u2026' sentence.\"\n}\n$rng.Collapse(0)   # wdCollapseEnd -> collapse to just after \"...e-mail address\"\n$rng.InsertAfter(\"es\")\n\n# 2) Insert the new e-mail address immediately before the existing one.\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$found2 = $rng2.Find.Execute(\"s.chopra@ucl.ac.uk\")\nif (-not $found2) {\n    throw \"Could not find 's.chopra@ucl.ac.uk' to anchor the new address on.\"\n}\n$rng2.InsertBefore(\"m.xochicale@ucl.ac.uk, \")\n\n# 3) Add the trailing space that now follows the last e-mail address.\n$rng3 = $d.Content\n$rng3.Find.ClearFormatting()\n$found3 = $rng3.Find.Execute(\"s.chopra@ucl.ac.uk\")\nif (-not $found3) {\n    throw \"Could not find 's.chopra@ucl.ac.uk' to append the trailing space after.\"\n}\n$rng3.Collapse(0)  # wdCollapseEnd -> collapse to just after the address\n$rng3.InsertAfter(\" \")\n"}
